$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.316.56'
$ws.Range('E2').Value = '  +1.03%  '

$ws.Range('D3').Value = '1.619.80'
$ws.Range('E3').Value = '  +1.64%  '

$ws.Range('E4').Value = '  -0.03%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '212.12'
$c.ClearFormats()
$ws.Range('E5').Value = '  +0.65%  '

$ws.Range('E6').Value = '  -0.09%  '

$ws.Range('E7').Value = '  +0.75%  '

$ws.Range('E8').Value = '  +0.27%  '

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.0616'
$c.ClearFormats()
$ws.Range('E9').Value = '  +0.30%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '18.85'
$c.ClearFormats()
$ws.Range('E10').Value = '  +4.56%  '

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0813'
$c.ClearFormats()
$ws.Range('E11').Value = '  +0.56%  '

$ws.Range('D12').Value = '1.844.24'
$ws.Range('E12').Value = '  +1.71%  '

$ws.Range('D13').Value = '1.623.59'
$ws.Range('E13').Value = '  +2.01%  '

$ws.Range('E14').Value = '  +0.52%  '

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.519'
$c.ClearFormats()
$ws.Range('E15').Value = '  +0.97%  '

$ws.Range('D16').Value = '26.308.80'
$ws.Range('E16').Value = '  +1.12%  '

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '62.51'
$c.ClearFormats()
$ws.Range('E17').Value = '  +3.53%  '

$ws.Range('D18').Value = '0.0₃0727'
$ws.Range('E18').Value = '  +0.32%  '

$ws.Range('E19').Value = '  -0.08%  '

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '201.96'
$c.ClearFormats()
$ws.Range('E20').Value = '  -0.14%  '

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '4.28'
$c.ClearFormats()
$ws.Range('E21').Value = '  +0.36%  '

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '9.34'
$c.ClearFormats()
$ws.Range('E22').Value = '  +0.85%  '

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '6.04'
$c.ClearFormats()
$ws.Range('E23').Value = '  +0.48%  '

$ws.Range('E24').Value = '  -4.79%  '

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '144.44'
$c.ClearFormats()
$ws.Range('E25').Value = '  +0.37%  '

$ws.Range('E26').Value = '  -0.08%  '

$ws.Range('E27').Value = '  -1.59%  '

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '15.18'
$c.ClearFormats()
$ws.Range('E28').Value = '  +0.26%  '

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '6.57'
$c.ClearFormats()
$ws.Range('E29').Value = '  +1.19%  '

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.0513'
$c.ClearFormats()
$ws.Range('E30').Value = '  +7.95%  '

$ws.Range('E31').Value = '  +0.44%  '

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.18'
$c.ClearFormats()
$ws.Range('E32').Value = '  +1.33%  '

$ws.Range('E33').Value = '  +0.00%  '

$ws.Range('E34').Value = '  +0.59%  '

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '2.39'
$c.ClearFormats()
$ws.Range('E35').Value = '  +2.21%  '

$ws.Range('D36').Value = '1.177.86'
$ws.Range('E36').Value = '  +4.01%  '

$ws.Range('E37').Value = '  +0.27%  '

$ws.Range('E38').Value = '  +1.69%  '

$ws.Range('E39').Value = '  +0.04%  '

$ws.Range('E40').Value = '  -0.03%  '

$ws.Range('E41').Value = '  +0.95%  '

$ws.Range('E42').Value = '  +4.18%  '

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.787'
$c.ClearFormats()
$ws.Range('E43').Value = '  +1.09%  '

$ws.Range('D44').Value = '1.756.46'
$ws.Range('E44').Value = '  +1.88%  '

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '92.56'
$c.ClearFormats()
$ws.Range('E45').Value = '  +0.42%  '

$ws.Range('E46').Value = '  +2.69%  '

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '53.86'
$c.ClearFormats()
$ws.Range('E47').Value = '  -0.31%  '

$ws.Range('E48').Value = '  +0.63%  '

$ws.Range('E49').Value = '  +0.34%  '

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E50').Value = '  -0.30%  '

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '7.30'
$c.ClearFormats()
$ws.Range('E51').Value = '  +2.63%  '

